$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (F column) counts
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 298
$ws1.Range("F3").Value = 1262
$ws1.Range("F4").Value = 2760

# Sheet "全部类型" (all types) - same rows mirrored here
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 298
$ws4.Range("F5").Value = 1262
$ws4.Range("F6").Value = 2760
